$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Add the new "TGmonitor7" sheet at the very end of the workbook (after the
# current last sheet, commonHW_AI) and make it the active sheet/tab.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "TGmonitor7"

# Column widths (approximate the source document's 42.22 / 45.61 char widths)
$ws.Columns.Item(1).ColumnWidth = 41.3
$ws.Columns.Item(2).ColumnWidth = 44.76

# ---------------------------------------------------------------------------
# Content - two-column key/value table describing the TGmonitor7 hardware.
# ---------------------------------------------------------------------------

# NAPAJENI (power) section
$ws.Cells.Item(1,1).Value = "NAPÁJENÍ"
$ws.Cells.Item(1,2).Value = " "
$ws.Cells.Item(2,1).Value = "Napájecí napětí"
$ws.Cells.Item(2,2).Value = "24 V DC (± 20 %)"
$ws.Cells.Item(3,1).Value = "Doporučený nap. zdroj"
$ws.Cells.Item(3,2).Value = "min. 300 mA"

# row 4 intentionally blank (section separator)

# KONEKTORY (connectors) section
$ws.Cells.Item(5,1).Value = "KONEKTORY"
$ws.Cells.Item(5,2).Value = " "
$ws.Cells.Item(6,1).Value = "USB"
$ws.Cells.Item(6,2).Value = "4 x USB 2.0, microUSB"
$ws.Cells.Item(7,1).Value = "HDMI"
$ws.Cells.Item(7,2).Value = "standard A"
$ws.Cells.Item(8,1).Value = "Napájení"
$ws.Cells.Item(8,2).Value = "1 x 4pin WEIDMÜLLER BLF 2.50/04/180 SN BK BX"

# row 9 intentionally blank (section separator)

# DISPLEJ (display) section
$ws.Cells.Item(10,1).Value = "DISPLEJ"
$ws.Cells.Item(10,2).Value = " "
$ws.Cells.Item(11,1).Value = "Uhlopříčka"
$ws.Cells.Item(11,2).Value = "7 palců"
$ws.Cells.Item(12,1).Value = "Pozoravací úhly"
$ws.Cells.Item(12,2).Value = "170 °"
$ws.Cells.Item(13,1).Value = "Rozlišení"
$ws.Cells.Item(13,2).Value = "1024 x 600 px"
$ws.Cells.Item(14,1).Value = "Rozměr zobrazovací plochy"
$ws.Cells.Item(14,2).Value = "154,21 x 85,92 mm"
$ws.Cells.Item(15,1).Value = "Rozteč pixelů"
$ws.Cells.Item(15,2).Value = "150,6(H) x 143,2(V) um"
$ws.Cells.Item(16,1).Value = "Barevný gamut"
$ws.Cells.Item(16,2).Value = "45% NTSC"
$ws.Cells.Item(17,1).Value = "Maximální jas"
$ws.Cells.Item(17,2).Value = "300 cd/m²"
$ws.Cells.Item(18,1).Value = "Kontrast"
$ws.Cells.Item(18,2).Value = "800:1"

# ---------------------------------------------------------------------------
# Formatting
# ---------------------------------------------------------------------------

# Rows 1-10: general number format, no wrap (same as the rest of the workbook)
$ws.Range("A1:B10").WrapText = $false

# Rows 11-23: wrapped text cells (A11:B23), including the trailing blank rows
# reserved for future chapters (dimensions, mounting, ...)
$ws.Range("A11:B23").WrapText = $true

# B16 ("45% NTSC") kept the donor cell's percentage-like number format
$ws.Range("B16").NumberFormat = "0.00\ %"

# B18 ("800:1") kept the donor cell's explicit text number format
$ws.Range("B18").NumberFormat = "@"

# Row heights match the rest of the workbook (12.8pt rows)
$ws.Range("A1:B23").EntireRow.RowHeight = 12.8

# Page setup (margins / paper size / header & footer) matching the rest of
# the workbook's sheets
$ws.PageSetup.LeftMargin = 56.7
$ws.PageSetup.RightMargin = 56.7
$ws.PageSetup.TopMargin = 75.8
$ws.PageSetup.BottomMargin = 75.8
$ws.PageSetup.HeaderMargin = 56.7
$ws.PageSetup.FooterMargin = 56.7
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.CenterHeader = '&"Times New Roman,obyčejné"&12&A'
$ws.PageSetup.CenterFooter = '&"Times New Roman,obyčejné"&12Stránka &P'

# Make the new sheet the active tab, matching the saved workbook view state
$ws.Activate()
$ws.Range("B21").Select()
